{"js": "/*\n * Applies the \"addition_and_subtraction_within_100\" worksheet update:\n *   - top date paragraph: \"2025-12-10 Wednesday\" -> \"2025-12-11 Thursday\"\n *   - each of the 100 arithmetic-expression cells (20 rows x 5 cols) gets\n *     its old expression text replaced with the new one, in document order.\n * Replacement is done in-place on each cell's Range.Text (via\n * TableCell.value), which preserves the existing run formatting\n * (rFonts/sz) exactly as in the target diff.\n */\nconst pairs = [[\"2025-12-10 Wednesday\", \"2025-12-11 Thursday\"], [\"78-33=\", \"90-43=\"], [\"5+52=\", \"7+80=\"], [\"6+43=\", \"41-4=\"], [\"10-2=\", \"2+85=\"], [\"41-3=\", \"55-29=\"], [\"98-88=\", \"73-18=\"], [\"41+56=\", \"68+24=\"], [\"70+17=\", \"44-36=\"], [\"50-25=\", \"38+44=\"], [\"51-46=\", \"46+22=\"], [\"98-5=\", \"13+76=\"], [\"55-30=\", \"66-2=\"], [\"73+0=\", \"77-15=\"], [\"31+54=\", \"89-24=\"], [\"78+10=\", \"46+6=\"], [\"43+9=\", \"74-2=\"], [\"99-85=\", \"3+59=\"], [\"79-34=\", \"37-20=\"], [\"3+87=\", \"38+54=\"], [\"84-81=\", \"23+16=\"], [\"87-22=\", \"35-14=\"], [\"98-16=\", \"73-67=\"], [\"28+5=\", \"8+45=\"], [\"33-8=\", \"6+48=\"], [\"59-40=\", \"46+4=\"], [\"36+49=\", \"50-0=\"], [\"67-36=\", \"52+44=\"], [\"31+59=\", \"66-53=\"], [\"41-16=\", \"73-28=\"], [\"75-61=\", \"86-57=\"], [\"8+35=\", \"58+14=\"], [\"34-24=\", \"84-21=\"], [\"4+54=\", \"31+35=\"], [\"72-18=\", \"14+22=\"], [\"58+31=\", \"82+3=\"], [\"4+52=\", \"51-24=\"], [\"27+64=\", \"10+76=\"], [\"9-4=\", \"43-36=\"], [\"74+11=\", \"36-7=\"], [\"58+41=\", \"27+25=\"], [\"37-21=\", \"18+6=\"], [\"17+31=\", \"69+22=\"], [\"65-64=\", \"22+21=\"], [\"10+81=\", \"15+39=\"], [\"62-19=\", \"76-39=\"], [\"47-10=\", \"16+27=\"], [\"83-11=\", \"6+40=\"], [\"99-36=\", \"54-37=\"], [\"24-24=\", \"74-20=\"], [\"48+12=\", \"33-17=\"], [\"24-20=\", \"7+17=\"], [\"63+28=\", \"43-24=\"], [\"98-3=\", \"45-41=\"], [\"66-39=\", \"6+55=\"], [\"36-19=\", \"47+17=\"], [\"81-23=\", \"97-63=\"], [\"82-6=\", \"2+38=\"], [\"57+39=\", \"76-69=\"], [\"3+49=\", \"65+9=\"], [\"41-31=\", \"16+20=\"], [\"4+35=\", \"93-39=\"], [\"90-67=\", \"76-25=\"], [\"52-43=\", \"8+40=\"], [\"34+34=\", \"71-59=\"], [\"3+84=\", \"78-63=\"], [\"44+10=\", \"91-77=\"], [\"80-33=\", \"18+55=\"], [\"39+53=\", \"3+25=\"], [\"53-40=\", \"7+63=\"], [\"88+4=\", \"95-83=\"], [\"41-35=\", \"86-84=\"], [\"50-43=\", \"54+40=\"], [\"57+10=\", \"65-49=\"], [\"84-62=\", \"39-6=\"], [\"63+2=\", \"78-43=\"], [\"96-38=\", \"24+54=\"], [\"24+5=\", \"66+23=\"], [\"54+39=\", \"44+41=\"], [\"36+34=\", \"68-62=\"], [\"17-1=\", \"66-17=\"], [\"77+19=\", \"10+31=\"], [\"11+19=\", \"7+1=\"], [\"97-89=\", \"54+28=\"], [\"36+5=\", \"9+56=\"], [\"60-41=\", \"65+6=\"], [\"90-60=\", \"98-87=\"], [\"3+52=\", \"45-34=\"], [\"65-9=\", \"17+34=\"], [\"25+13=\", \"48-38=\"], [\"35+57=\", \"44-5=\"], [\"70-42=\", \"32-8=\"], [\"53-2=\", \"5+50=\"], [\"89-76=\", \"55-32=\"], [\"4+42=\", \"40+2=\"], [\"15-13=\", \"90-1=\"], [\"53+23=\", \"89-47=\"], [\"33+58=\", \"36+32=\"], [\"4+67=\", \"48+24=\"], [\"59-21=\", \"36-20=\"], [\"64-61=\", \"49-44=\"]];\nconst datePair = pairs[0];\nconst cellPairs = pairs.slice(1);\nconst ROWS = 20, COLS = 5;\n\n// --- 1. Update the date paragraph (first paragraph in the body) ---\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst dateParagraph = paragraphs.items[0];\ndateParagraph.load(\"text\");\nawait context.sync();\n\nif (dateParagraph.text === datePair[0]) {\n  dateParagraph.insertText(datePair[1], \"Replace\");\n} else {\n  // Fall back to a text search if the first paragraph doesn't match exactly.\n  const hits = body.search(datePair[0], { matchCase: true });\n  hits.load(\"items\");\n  await context.sync();\n  if (hits.items.length > 0) {\n    hits.items[0].insertText(datePair[1], \"Replace\");\n  }\n}\nawait context.sync();\n\n// --- 2. Update every table cell in document order ---\nconst tables = body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst cells = [];\nfor (let r = 0; r < ROWS; r++) {\n  for (let c = 0; c < COLS; c++) {\n    const cell = table.getCell(r, c);\n    cell.load(\"value\");\n    cells.push(cell);\n  }\n}\nawait context.sync();\n\nfor (let i = 0; i < cellPairs.length; i++) {\n  const [oldText, newText] = cellPairs[i];\n  const cell = cells[i];\n  if (cell.value === oldText) {\n    cell.value = newText;\n  } else {\n    // Defensive fallback: search within the cell's range for the old text.\n    const cellRange = cell.getRange();\n    const hits = cellRange.search(oldText, { matchCase: true });\n    hits.load(\"items\");\n    await context.sync();\n    if (hits.items.length > 0) {\n      hits.items[0].insertText(newText, \"Replace\");\n    }\n  }\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Ordered list of (old, new) text replacements, taken from document order:\n# 1 date heading + 100 arithmetic-expression table cells.\n$pairs = @(\n  @('2025-12-10 Wednesday', '2025-12-11 Thursday'),\n  @('78-33=', '90-43='),\n  @('5+52=', '7+80='),\n  @('6+43=', '41-4='),\n  @('10-2=', '2+85='),\n  @('41-3=', '55-29='),\n  @('98-88=', '73-18='),\n  @('41+56=', '68+24='),\n  @('70+17=', '44-36='),\n  @('50-25=', '38+44='),\n  @('51-46=', '46+22='),\n  @('98-5=', '13+76='),\n  @('55-30=', '66-2='),\n  @('73+0=', '77-15='),\n  @('31+54=', '89-24='),\n  @('78+10=', '46+6='),\n  @('43+9=', '74-2='),\n  @('99-85=', '3+59='),\n  @('79-34=', '37-20='),\n  @('3+87=', '38+54='),\n  @('84-81=', '23+16='),\n  @('87-22=', '35-14='),\n  @('98-16=', '73-67='),\n  @('28+5=', '8+45='),\n  @('33-8=', '6+48='),\n  @('59-40=', '46+4='),\n  @('36+49=', '50-0='),\n  @('67-36=', '52+44='),\n  @('31+59=', '66-53='),\n  @('41-16=', '73-28='),\n  @('75-61=', '86-57='),\n  @('8+35=', '58+14='),\n  @('34-24=', '84-21='),\n  @('4+54=', '31+35='),\n  @('72-18=', '14+22='),\n  @('58+31=', '82+3='),\n  @('4+52=', '51-24='),\n  @('27+64=', '10+76='),\n  @('9-4=', '43-36='),\n  @('74+11=', '36-7='),\n  @('58+41=', '27+25='),\n  @('37-21=', '18+6='),\n  @('17+31=', '69+22='),\n  @('65-64=', '22+21='),\n  @('10+81=', '15+39='),\n  @('62-19=', '76-39='),\n  @('47-10=', '16+27='),\n  @('83-11=', '6+40='),\n  @('99-36=', '54-37='),\n  @('24-24=', '74-20='),\n  @('48+12=', '33-17='),\n  @('24-20=', '7+17='),\n  @('63+28=', '43-24='),\n  @('98-3=', '45-41='),\n  @('66-39=', '6+55='),\n  @('36-19=', '47+17='),\n  @('81-23=', '97-63='),\n  @('82-6=', '2+38='),\n  @('57+39=', '76-69='),\n  @('3+49=', '65+9='),\n  @('41-31=', '16+20='),\n  @('4+35=', '93-39='),\n  @('90-67=', '76-25='),\n  @('52-43=', '8+40='),\n  @('34+34=', '71-59='),\n  @('3+84=', '78-63='),\n  @('44+10=', '91-77='),\n  @('80-33=', '18+55='),\n  @('39+53=', '3+25='),\n  @('53-40=', '7+63='),\n  @('88+4=', '95-83='),\n  @('41-35=', '86-84='),\n  @('50-43=', '54+40='),\n  @('57+10=', '65-49='),\n  @('84-62=', '39-6='),\n  @('63+2=', '78-43='),\n  @('96-38=', '24+54='),\n  @('24+5=', '66+23='),\n  @('54+39=', '44+41='),\n  @('36+34=', '68-62='),\n  @('17-1=', '66-17='),\n  @('77+19=', '10+31='),\n  @('11+19=', '7+1='),\n  @('97-89=', '54+28='),\n  @('36+5=', '9+56='),\n  @('60-41=', '65+6='),\n  @('90-60=', '98-87='),\n  @('3+52=', '45-34='),\n  @('65-9=', '17+34='),\n  @('25+13=', '48-38='),\n  @('35+57=', '44-5='),\n  @('70-42=', '32-8='),\n  @('53-2=', '5+50='),\n  @('89-76=', '55-32='),\n  @('4+42=', '40+2='),\n  @('15-13=', '90-1='),\n  @('53+23=', '89-47='),\n  @('33+58=', '36+32='),\n  @('4+67=', '48+24='),\n  @('59-21=', '36-20='),\n  @('64-61=', '49-44=')\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $range = $d.Content\n  [void]$range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
